$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 19 de Abril de 2020 a las 12:52"

# Update province/city case data (values refreshed + re-sorted descending by "Casos totales")
$ws.Cells.Item(4, 1).Value = "Madrid"
$ws.Cells.Item(4, 2).Value = 54884
$ws.Cells.Item(4, 3).Value = 31313
$ws.Cells.Item(4, 4).Value = 16332
$ws.Cells.Item(4, 5).Value = 7239
$ws.Cells.Item(5, 1).Value = "Cataluña"
$ws.Cells.Item(5, 2).Value = 40988
$ws.Cells.Item(5, 3).Value = 19088
$ws.Cells.Item(5, 4).Value = 17777
$ws.Cells.Item(5, 5).Value = 4123
$ws.Cells.Item(6, 1).Value = "Galicia"
$ws.Cells.Item(6, 2).Value = 8125
$ws.Cells.Item(6, 3).Value = 1536
$ws.Cells.Item(6, 4).Value = 6145
$ws.Cells.Item(6, 5).Value = 444
$ws.Cells.Item(7, 1).Value = "Bizkaia/Vizcaya"
$ws.Cells.Item(7, 2).Value = 6667
$ws.Cells.Item(7, 3).Value = 6144
$ws.Cells.Item(7, 4).Value = 4953
$ws.Cells.Item(7, 5).Value = 515
$ws.Cells.Item(8, 1).Value = "Ciudad Real"
$ws.Cells.Item(8, 2).Value = 6300
$ws.Cells.Item(8, 3).Value = 3838
$ws.Cells.Item(8, 4).Value = 10545
$ws.Cells.Item(8, 5).Value = 770
$ws.Cells.Item(9, 1).Value = "Valencia/Valencia"
$ws.Cells.Item(9, 2).Value = 5131
$ws.Cells.Item(9, 3).Value = 2194
$ws.Cells.Item(9, 4).Value = 2583
$ws.Cells.Item(9, 5).Value = 515
$ws.Cells.Item(10, 1).Value = "Navarra"
$ws.Cells.Item(10, 2).Value = 4621
$ws.Cells.Item(10, 3).Value = 992
$ws.Cells.Item(10, 4).Value = 3260
$ws.Cells.Item(10, 5).Value = 369
$ws.Cells.Item(11, 1).Value = "Toledo"
$ws.Cells.Item(11, 2).Value = 3908
$ws.Cells.Item(11, 3).Value = 3838
$ws.Cells.Item(11, 4).Value = 10545
$ws.Cells.Item(11, 5).Value = 497
$ws.Cells.Item(12, 1).Value = "Albacete"
$ws.Cells.Item(12, 2).Value = 3732
$ws.Cells.Item(12, 3).Value = 3838
$ws.Cells.Item(12, 4).Value = 10545
$ws.Cells.Item(12, 5).Value = 358
$ws.Cells.Item(13, 1).Value = "La Rioja"
$ws.Cells.Item(13, 2).Value = 3705
$ws.Cells.Item(13, 3).Value = 1612
$ws.Cells.Item(13, 4).Value = 1175
$ws.Cells.Item(13, 5).Value = 277
$ws.Cells.Item(14, 1).Value = "Zaragoza"
$ws.Cells.Item(14, 2).Value = 3643
$ws.Cells.Item(14, 3).Value = 911
$ws.Cells.Item(14, 4).Value = 2258
$ws.Cells.Item(14, 5).Value = 474
$ws.Cells.Item(15, 1).Value = "Alacant/Alicante"
$ws.Cells.Item(15, 2).Value = 3476
$ws.Cells.Item(15, 3).Value = 1677
$ws.Cells.Item(15, 4).Value = 1443
$ws.Cells.Item(15, 5).Value = 388
$ws.Cells.Item(16, 1).Value = "Araba/Alava"
$ws.Cells.Item(16, 2).Value = 3156
$ws.Cells.Item(16, 3).Value = 6144
$ws.Cells.Item(16, 4).Value = 4953
$ws.Cells.Item(16, 5).Value = 308
$ws.Cells.Item(17, 1).Value = "Valladolid"
$ws.Cells.Item(17, 2).Value = 3098
$ws.Cells.Item(17, 3).Value = 1045
$ws.Cells.Item(17, 4).Value = 1813
$ws.Cells.Item(17, 5).Value = 243
$ws.Cells.Item(18, 1).Value = "Castilla-La Mancha"
$ws.Cells.Item(18, 2).Value = 2780
$ws.Cells.Item(18, 3).Value = 71
$ws.Cells.Item(18, 4).Value = 2446
$ws.Cells.Item(18, 5).Value = 263
$ws.Cells.Item(19, 1).Value = "Salamanca"
$ws.Cells.Item(19, 2).Value = 2514
$ws.Cells.Item(19, 3).Value = 761
$ws.Cells.Item(19, 4).Value = 1431
$ws.Cells.Item(19, 5).Value = 276
$ws.Cells.Item(20, 1).Value = "Malaga"
$ws.Cells.Item(20, 2).Value = 2472
$ws.Cells.Item(20, 3).Value = 835
$ws.Cells.Item(20, 4).Value = 1420
$ws.Cells.Item(20, 5).Value = 217
$ws.Cells.Item(21, 1).Value = "Sevilla"
$ws.Cells.Item(21, 2).Value = 2299
$ws.Cells.Item(21, 3).Value = 432
$ws.Cells.Item(21, 4).Value = 1660
$ws.Cells.Item(21, 5).Value = 207
$ws.Cells.Item(22, 1).Value = "Asturias"
$ws.Cells.Item(22, 2).Value = 2298
$ws.Cells.Item(22, 3).Value = 596
$ws.Cells.Item(22, 4).Value = 1506
$ws.Cells.Item(22, 5).Value = 196
$ws.Cells.Item(23, 1).Value = "Segovia"
$ws.Cells.Item(23, 2).Value = 2285
$ws.Cells.Item(23, 3).Value = 636
$ws.Cells.Item(23, 4).Value = 1469
$ws.Cells.Item(23, 5).Value = 167
$ws.Cells.Item(24, 1).Value = "Leon"
$ws.Cells.Item(24, 2).Value = 2285
$ws.Cells.Item(24, 3).Value = 1031
$ws.Cells.Item(24, 4).Value = 927
$ws.Cells.Item(24, 5).Value = 290
$ws.Cells.Item(25, 1).Value = "Gipuzkoa/Guipuzcoa"
$ws.Cells.Item(25, 2).Value = 2266
$ws.Cells.Item(25, 3).Value = 6144
$ws.Cells.Item(25, 4).Value = 4953
$ws.Cells.Item(25, 5).Value = 197
$ws.Cells.Item(26, 1).Value = "Caceres"
$ws.Cells.Item(26, 2).Value = 2067
$ws.Cells.Item(26, 3).Value = 340
$ws.Cells.Item(26, 4).Value = 1503
$ws.Cells.Item(26, 5).Value = 316
$ws.Cells.Item(27, 1).Value = "Cantabria"
$ws.Cells.Item(27, 2).Value = 2050
$ws.Cells.Item(27, 3).Value = 610
$ws.Cells.Item(27, 4).Value = 1257
$ws.Cells.Item(27, 5).Value = 153
$ws.Cells.Item(28, 1).Value = "A Coruña"
$ws.Cells.Item(28, 2).Value = 1969
$ws.Cells.Item(28, 3).Value = 333
$ws.Cells.Item(28, 4).Value = 1788
$ws.Cells.Item(28, 5).Value = 67
$ws.Cells.Item(29, 1).Value = "Granada"
$ws.Cells.Item(29, 2).Value = 1969
$ws.Cells.Item(29, 3).Value = 563
$ws.Cells.Item(29, 4).Value = 1213
$ws.Cells.Item(29, 5).Value = 193
$ws.Cells.Item(30, 1).Value = "Murcia"
$ws.Cells.Item(30, 2).Value = 1644
$ws.Cells.Item(30, 3).Value = 652
$ws.Cells.Item(30, 4).Value = 876
$ws.Cells.Item(30, 5).Value = 116
$ws.Cells.Item(31, 1).Value = "Pontevedra"
$ws.Cells.Item(31, 2).Value = 1536
$ws.Cells.Item(31, 3).Value = 333
$ws.Cells.Item(31, 4).Value = 1411
$ws.Cells.Item(31, 5).Value = 30
$ws.Cells.Item(32, 1).Value = "Burgos"
$ws.Cells.Item(32, 2).Value = 1522
$ws.Cells.Item(32, 3).Value = 620
$ws.Cells.Item(32, 4).Value = 685
$ws.Cells.Item(32, 5).Value = 160
$ws.Cells.Item(33, 1).Value = "Guadalajara"
$ws.Cells.Item(33, 2).Value = 1400
$ws.Cells.Item(33, 3).Value = 3838
$ws.Cells.Item(33, 4).Value = 10545
$ws.Cells.Item(33, 5).Value = 184
$ws.Cells.Item(34, 1).Value = "Tenerife"
$ws.Cells.Item(34, 2).Value = 1378
$ws.Cells.Item(34, 3).Value = 489
$ws.Cells.Item(34, 4).Value = 808
$ws.Cells.Item(34, 5).Value = 81
$ws.Cells.Item(35, 1).Value = "Cuenca"
$ws.Cells.Item(35, 2).Value = 1285
$ws.Cells.Item(35, 3).Value = 3838
$ws.Cells.Item(35, 4).Value = 10545
$ws.Cells.Item(35, 5).Value = 154
$ws.Cells.Item(36, 1).Value = "Jaen"
$ws.Cells.Item(36, 2).Value = 1274
$ws.Cells.Item(36, 3).Value = 252
$ws.Cells.Item(36, 4).Value = 891
$ws.Cells.Item(36, 5).Value = 131
$ws.Cells.Item(37, 1).Value = "Cordoba"
$ws.Cells.Item(37, 2).Value = 1266
$ws.Cells.Item(37, 3).Value = 353
$ws.Cells.Item(37, 4).Value = 841
$ws.Cells.Item(37, 5).Value = 72
$ws.Cells.Item(38, 1).Value = "Castello/Castellon"
$ws.Cells.Item(38, 2).Value = 1257
$ws.Cells.Item(38, 3).Value = 435
$ws.Cells.Item(38, 4).Value = 739
$ws.Cells.Item(38, 5).Value = 139
$ws.Cells.Item(39, 1).Value = "Soria"
$ws.Cells.Item(39, 2).Value = 1231
$ws.Cells.Item(39, 3).Value = 293
$ws.Cells.Item(39, 4).Value = 802
$ws.Cells.Item(39, 5).Value = 94
$ws.Cells.Item(40, 1).Value = "Cadiz"
$ws.Cells.Item(40, 2).Value = 1122
$ws.Cells.Item(40, 3).Value = 280
$ws.Cells.Item(40, 4).Value = 768
$ws.Cells.Item(40, 5).Value = 74
$ws.Cells.Item(41, 1).Value = "Avila"
$ws.Cells.Item(41, 2).Value = 1090
$ws.Cells.Item(41, 3).Value = 446
$ws.Cells.Item(41, 4).Value = 531
$ws.Cells.Item(41, 5).Value = 107
$ws.Cells.Item(42, 1).Value = "Badajoz"
$ws.Cells.Item(42, 2).Value = 972
$ws.Cells.Item(42, 3).Value = 396
$ws.Cells.Item(42, 4).Value = 504
$ws.Cells.Item(42, 5).Value = 72
$ws.Cells.Item(43, 1).Value = "Aragon"
$ws.Cells.Item(43, 2).Value = 907
$ws.Cells.Item(43, 3).Value = 29
$ws.Cells.Item(43, 4).Value = 838
$ws.Cells.Item(43, 5).Value = 40
$ws.Cells.Item(44, 1).Value = "Ourense"
$ws.Cells.Item(44, 2).Value = 751
$ws.Cells.Item(44, 3).Value = 333
$ws.Cells.Item(44, 4).Value = 660
$ws.Cells.Item(44, 5).Value = 22
$ws.Cells.Item(45, 1).Value = "Palencia"
$ws.Cells.Item(45, 2).Value = 673
$ws.Cells.Item(45, 3).Value = 220
$ws.Cells.Item(45, 4).Value = 399
$ws.Cells.Item(45, 5).Value = 58
$ws.Cells.Item(46, 1).Value = "Zamora"
$ws.Cells.Item(46, 2).Value = 595
$ws.Cells.Item(46, 3).Value = 210
$ws.Cells.Item(46, 4).Value = 314
$ws.Cells.Item(46, 5).Value = 63
$ws.Cells.Item(47, 1).Value = "Lugo"
$ws.Cells.Item(47, 2).Value = 586
$ws.Cells.Item(47, 3).Value = 333
$ws.Cells.Item(47, 4).Value = 520
$ws.Cells.Item(47, 5).Value = 11
$ws.Cells.Item(48, 1).Value = "Huesca"
$ws.Cells.Item(48, 2).Value = 585
$ws.Cells.Item(48, 3).Value = 135
$ws.Cells.Item(48, 4).Value = 370
$ws.Cells.Item(48, 5).Value = 80
$ws.Cells.Item(49, 1).Value = "Teruel"
$ws.Cells.Item(49, 2).Value = 540
$ws.Cells.Item(49, 3).Value = 118
$ws.Cells.Item(49, 4).Value = 358
$ws.Cells.Item(49, 5).Value = 64
$ws.Cells.Item(50, 1).Value = "Gran Canaria"
$ws.Cells.Item(50, 2).Value = 496
$ws.Cells.Item(50, 3).Value = 235
$ws.Cells.Item(50, 4).Value = 228
$ws.Cells.Item(50, 5).Value = 33
$ws.Cells.Item(51, 1).Value = "Almeria"
$ws.Cells.Item(51, 2).Value = 444
$ws.Cells.Item(51, 3).Value = 126
$ws.Cells.Item(51, 4).Value = 277
$ws.Cells.Item(51, 5).Value = 41
$ws.Cells.Item(52, 1).Value = "Huelva"
$ws.Cells.Item(52, 2).Value = 358
$ws.Cells.Item(52, 3).Value = 107
$ws.Cells.Item(52, 4).Value = 219
$ws.Cells.Item(52, 5).Value = 32
$ws.Cells.Item(53, 1).Value = "Mallorca"
$ws.Cells.Item(53, 2).Value = 210
$ws.Cells.Item(53, 3).Value = 18
$ws.Cells.Item(53, 4).Value = 194
$ws.Cells.Item(53, 5).Value = 12
$ws.Cells.Item(54, 1).Value = "Ceuta"
$ws.Cells.Item(54, 2).Value = 109
$ws.Cells.Item(54, 3).Value = 59
$ws.Cells.Item(54, 4).Value = 46
$ws.Cells.Item(54, 5).Value = 4
$ws.Cells.Item(55, 1).Value = "Melilla"
$ws.Cells.Item(55, 2).Value = 104
$ws.Cells.Item(55, 3).Value = 44
$ws.Cells.Item(55, 4).Value = 58
$ws.Cells.Item(55, 5).Value = 2
$ws.Cells.Item(56, 1).Value = "La Palma"
$ws.Cells.Item(56, 2).Value = 73
$ws.Cells.Item(56, 3).Value = 23
$ws.Cells.Item(56, 4).Value = 47
$ws.Cells.Item(56, 5).Value = 3
$ws.Cells.Item(57, 1).Value = "Lanzarote"
$ws.Cells.Item(57, 2).Value = 68
$ws.Cells.Item(57, 3).Value = 18
$ws.Cells.Item(57, 4).Value = 48
$ws.Cells.Item(57, 5).Value = 2
$ws.Cells.Item(58, 1).Value = "Igualada, Vilanova del Cami, Santa Margarida de Montbui y Odena"
$ws.Cells.Item(58, 2).Value = 58
$ws.Cells.Item(58, 3).Value = 0
$ws.Cells.Item(58, 4).Value = 58
$ws.Cells.Item(58, 5).Value = 3
$ws.Cells.Item(59, 1).Value = "Fuerteventura"
$ws.Cells.Item(59, 2).Value = 24
$ws.Cells.Item(59, 3).Value = 18
$ws.Cells.Item(59, 4).Value = 6
$ws.Cells.Item(59, 5).Value = 0
$ws.Cells.Item(60, 1).Value = "Ibiza"
$ws.Cells.Item(60, 2).Value = 21
$ws.Cells.Item(60, 3).Value = 18
$ws.Cells.Item(60, 4).Value = 20
$ws.Cells.Item(60, 5).Value = 1
$ws.Cells.Item(61, 1).Value = "Menorca"
$ws.Cells.Item(61, 2).Value = 15
$ws.Cells.Item(61, 3).Value = 18
$ws.Cells.Item(61, 4).Value = 13
$ws.Cells.Item(61, 5).Value = 0
$ws.Cells.Item(62, 1).Value = "Arroyo de la Luz"
$ws.Cells.Item(62, 2).Value = 7
$ws.Cells.Item(62, 3).Value = 0
$ws.Cells.Item(62, 4).Value = 7
$ws.Cells.Item(62, 5).Value = 0
$ws.Cells.Item(63, 1).Value = "La Gomera"
$ws.Cells.Item(63, 2).Value = 7
$ws.Cells.Item(63, 3).Value = 5
$ws.Cells.Item(63, 4).Value = 2
$ws.Cells.Item(63, 5).Value = 0
$ws.Cells.Item(64, 1).Value = "El Hierro"
$ws.Cells.Item(64, 2).Value = 1
$ws.Cells.Item(64, 3).Value = 1
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(64, 5).Value = 0
$ws.Cells.Item(65, 1).Value = "Formentera"
$ws.Cells.Item(65, 2).Value = 0
$ws.Cells.Item(65, 3).Value = 10
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(65, 5).Value = 8
